$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 104: correct the date/time value in column A ---
$ws.Range("A104").Value = 45491.2916666667

# --- Append new row 105 ---

# Column A is a date/time; clone A104's formatting (style 1) onto A105 first
# so we don't mint a duplicate style entry, then overwrite with the new value.
$ws.Range("A104").Copy($ws.Range("A105"))
$ws.Range("A105").Value = 45492.3665509259

$ws.Range("B105").Value = 1500
$ws.Range("C105").Value = 3.19000005722046
$ws.Range("D105").Value = 3.11999988555908
$ws.Range("E105").Value = 3.11999988555908
$ws.Range("F105").Value = 3.19000005722046

# Column G stores a textual representation of the close price (shared string,
# not a plain number, matching the rest of this column). Force text typing by
# switching to a text number format while assigning the value, then restore
# the plain "Normal" style so G105 ends up with the default style like the
# other string cells in this column.
$ws.Range("G105").NumberFormat = "@"
$ws.Range("G105").Value = "3.19000005722046"
$ws.Range("G105").Style = "Normal"

$ws.Range("H105").Value = "ESPE.MI"
